$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 26
$ws.Cells.Item(26, 2).Value = 6802962
$ws.Cells.Item(26, 5).Value = 'FC Trinity Zlin'
$ws.Cells.Item(26, 6).Value = 'Slovan Liberec'
$ws.Cells.Item(26, 7).Value = 1
$ws.Cells.Item(26, 8).Value = 1
$ws.Cells.Item(26, 9).Value = 1
$ws.Cells.Item(26, 10).Value = 1
$ws.Cells.Item(26, 11).Value = 'D'
$ws.Cells.Item(26, 12).Value = 2.6
$ws.Cells.Item(26, 13).Value = 3.2
$ws.Cells.Item(26, 14).Value = 2.4
$ws.Cells.Item(26, 15).Value = 2.75
$ws.Cells.Item(26, 16).Value = 3.25
$ws.Cells.Item(26, 17).Value = 2.25
$ws.Cells.Item(26, 18).Value = 0.25
$ws.Cells.Item(26, 19).Value = 1.8
$ws.Cells.Item(26, 20).Value = 2.05
$ws.Cells.Item(26, 21).Value = 2.5
$ws.Cells.Item(26, 22).Value = 1.825
$ws.Cells.Item(26, 23).Value = 2.025
$ws.Cells.Item(26, 24).Value = -1
$ws.Cells.Item(26, 25).Value = 2.25
$ws.Cells.Item(26, 27).Value = 0.4
$ws.Cells.Item(26, 28).Value = -0.5
$ws.Cells.Item(26, 30).Value = 1.025

# Row 27
$ws.Cells.Item(27, 2).Value = 6802958
$ws.Cells.Item(27, 5).Value = 'Banik Ostrava'
$ws.Cells.Item(27, 6).Value = 'Hradec Kralove'
$ws.Cells.Item(27, 7).Value = 2
$ws.Cells.Item(27, 8).Value = 0
$ws.Cells.Item(27, 9).Value = 0
$ws.Cells.Item(27, 10).Value = 0
$ws.Cells.Item(27, 11).Value = 'H'
$ws.Cells.Item(27, 12).Value = 1.666
$ws.Cells.Item(27, 13).Value = 3.5
$ws.Cells.Item(27, 14).Value = 4.333
$ws.Cells.Item(27, 15).Value = 1.75
$ws.Cells.Item(27, 16).Value = 3.75
$ws.Cells.Item(27, 17).Value = 3.8
$ws.Cells.Item(27, 18).Value = -0.5
$ws.Cells.Item(27, 19).Value = 1.825
$ws.Cells.Item(27, 20).Value = 2.025
$ws.Cells.Item(27, 21).Value = 2.75
$ws.Cells.Item(27, 22).Value = 1.925
$ws.Cells.Item(27, 23).Value = 1.925
$ws.Cells.Item(27, 24).Value = 0.75
$ws.Cells.Item(27, 25).Value = -1
$ws.Cells.Item(27, 27).Value = 0.825
$ws.Cells.Item(27, 28).Value = -1
$ws.Cells.Item(27, 30).Value = 0.925

# Row 46
$ws.Cells.Item(46, 2).Value = 7098155
$ws.Cells.Item(46, 5).Value = 'FK Jablonec'
$ws.Cells.Item(46, 6).Value = 'Slavia Prague'
$ws.Cells.Item(46, 7).Value = 1
$ws.Cells.Item(46, 8).Value = 1
$ws.Cells.Item(46, 9).Value = 1
$ws.Cells.Item(46, 11).Value = 'D'
$ws.Cells.Item(46, 12).Value = 7
$ws.Cells.Item(46, 13).Value = 4.5
$ws.Cells.Item(46, 14).Value = 1.363
$ws.Cells.Item(46, 15).Value = 7
$ws.Cells.Item(46, 16).Value = 4.5
$ws.Cells.Item(46, 17).Value = 1.4
$ws.Cells.Item(46, 18).Value = 1.25
$ws.Cells.Item(46, 19).Value = 1.95
$ws.Cells.Item(46, 20).Value = 1.9
$ws.Cells.Item(46, 21).Value = 3
$ws.Cells.Item(46, 24).Value = -1
$ws.Cells.Item(46, 25).Value = 3.5
$ws.Cells.Item(46, 27).Value = 0.95

# Row 47
$ws.Cells.Item(47, 2).Value = 6802978
$ws.Cells.Item(47, 5).Value = 'Sigma Olomouc'
$ws.Cells.Item(47, 6).Value = 'Slovan Liberec'
$ws.Cells.Item(47, 7).Value = 2
$ws.Cells.Item(47, 8).Value = 0
$ws.Cells.Item(47, 9).Value = 0
$ws.Cells.Item(47, 11).Value = 'H'
$ws.Cells.Item(47, 12).Value = 1.85
$ws.Cells.Item(47, 13).Value = 3.4
$ws.Cells.Item(47, 14).Value = 3.8
$ws.Cells.Item(47, 15).Value = 1.75
$ws.Cells.Item(47, 16).Value = 3.5
$ws.Cells.Item(47, 17).Value = 4.2
$ws.Cells.Item(47, 18).Value = -0.75
$ws.Cells.Item(47, 19).Value = 2.05
$ws.Cells.Item(47, 20).Value = 1.8
$ws.Cells.Item(47, 21).Value = 2.75
$ws.Cells.Item(47, 24).Value = 0.75
$ws.Cells.Item(47, 25).Value = -1
$ws.Cells.Item(47, 27).Value = 1.05

# Row 53
$ws.Cells.Item(53, 2).Value = 6802981
$ws.Cells.Item(53, 5).Value = 'Viktoria Plzen'
$ws.Cells.Item(53, 6).Value = 'Bohemians 1905'
$ws.Cells.Item(53, 7).Value = 2
$ws.Cells.Item(53, 8).Value = 0
$ws.Cells.Item(53, 9).Value = 2
$ws.Cells.Item(53, 11).Value = 'H'
$ws.Cells.Item(53, 12).Value = 1.45
$ws.Cells.Item(53, 13).Value = 4.333
$ws.Cells.Item(53, 14).Value = 5.5
$ws.Cells.Item(53, 15).Value = 1.3
$ws.Cells.Item(53, 16).Value = 5
$ws.Cells.Item(53, 17).Value = 7
$ws.Cells.Item(53, 18).Value = -1.5
$ws.Cells.Item(53, 19).Value = 1.925
$ws.Cells.Item(53, 20).Value = 1.925
$ws.Cells.Item(53, 21).Value = 3.25
$ws.Cells.Item(53, 22).Value = 2.025
$ws.Cells.Item(53, 23).Value = 1.825
$ws.Cells.Item(53, 24).Value = 0.3
$ws.Cells.Item(53, 26).Value = -1
$ws.Cells.Item(53, 27).Value = 0.925
$ws.Cells.Item(53, 28).Value = -1
$ws.Cells.Item(53, 30).Value = 0.825

# Row 54
$ws.Cells.Item(54, 2).Value = 6802987
$ws.Cells.Item(54, 5).Value = 'Slovan Liberec'
$ws.Cells.Item(54, 6).Value = 'Sparta Prague'
$ws.Cells.Item(54, 7).Value = 0
$ws.Cells.Item(54, 8).Value = 2
$ws.Cells.Item(54, 9).Value = 0
$ws.Cells.Item(54, 11).Value = 'A'
$ws.Cells.Item(54, 12).Value = 4.5
$ws.Cells.Item(54, 13).Value = 3.8
$ws.Cells.Item(54, 14).Value = 1.65
$ws.Cells.Item(54, 15).Value = 5.5
$ws.Cells.Item(54, 16).Value = 4
$ws.Cells.Item(54, 17).Value = 1.533
$ws.Cells.Item(54, 18).Value = 1
$ws.Cells.Item(54, 19).Value = 1.9
$ws.Cells.Item(54, 20).Value = 1.95
$ws.Cells.Item(54, 21).Value = 2.75
$ws.Cells.Item(54, 22).Value = 1.925
$ws.Cells.Item(54, 23).Value = 1.925
$ws.Cells.Item(54, 24).Value = -1
$ws.Cells.Item(54, 26).Value = 0.5329999999999999
$ws.Cells.Item(54, 27).Value = -1
$ws.Cells.Item(54, 28).Value = 0.95
$ws.Cells.Item(54, 30).Value = 0.925

# Row 59
$ws.Cells.Item(59, 2).Value = 6802993
$ws.Cells.Item(59, 5).Value = 'MFK Karvina'
$ws.Cells.Item(59, 6).Value = 'FK Teplice'
$ws.Cells.Item(59, 7).Value = 0
$ws.Cells.Item(59, 11).Value = 'A'
$ws.Cells.Item(59, 12).Value = 2.5
$ws.Cells.Item(59, 13).Value = 3.3
$ws.Cells.Item(59, 14).Value = 2.7
$ws.Cells.Item(59, 15).Value = 2.15
$ws.Cells.Item(59, 16).Value = 3.5
$ws.Cells.Item(59, 17).Value = 3.2
$ws.Cells.Item(59, 19).Value = 1.875
$ws.Cells.Item(59, 20).Value = 1.975
$ws.Cells.Item(59, 21).Value = 2.5
$ws.Cells.Item(59, 22).Value = 1.875
$ws.Cells.Item(59, 23).Value = 1.975
$ws.Cells.Item(59, 25).Value = -1
$ws.Cells.Item(59, 26).Value = 2.2
$ws.Cells.Item(59, 27).Value = -1
$ws.Cells.Item(59, 28).Value = 0.9750000000000001
$ws.Cells.Item(59, 30).Value = 0.9750000000000001

# Row 60
$ws.Cells.Item(60, 2).Value = 6802992
$ws.Cells.Item(60, 5).Value = 'Bohemians 1905'
$ws.Cells.Item(60, 6).Value = 'Banik Ostrava'
$ws.Cells.Item(60, 7).Value = 1
$ws.Cells.Item(60, 11).Value = 'D'
$ws.Cells.Item(60, 12).Value = 2.25
$ws.Cells.Item(60, 13).Value = 3.4
$ws.Cells.Item(60, 14).Value = 3
$ws.Cells.Item(60, 15).Value = 2.3
$ws.Cells.Item(60, 16).Value = 3.4
$ws.Cells.Item(60, 17).Value = 3
$ws.Cells.Item(60, 19).Value = 2
$ws.Cells.Item(60, 20).Value = 1.85
$ws.Cells.Item(60, 21).Value = 2.75
$ws.Cells.Item(60, 22).Value = 1.975
$ws.Cells.Item(60, 23).Value = 1.875
$ws.Cells.Item(60, 25).Value = 2.4
$ws.Cells.Item(60, 26).Value = -1
$ws.Cells.Item(60, 27).Value = -0.5
$ws.Cells.Item(60, 28).Value = 0.425
$ws.Cells.Item(60, 30).Value = 0.875

# Row 73
$ws.Cells.Item(73, 2).Value = 6803005
$ws.Cells.Item(73, 5).Value = 'Sigma Olomouc'
$ws.Cells.Item(73, 6).Value = 'Banik Ostrava'
$ws.Cells.Item(73, 7).Value = 0
$ws.Cells.Item(73, 8).Value = 3
$ws.Cells.Item(73, 9).Value = 0
$ws.Cells.Item(73, 10).Value = 2
$ws.Cells.Item(73, 11).Value = 'A'
$ws.Cells.Item(73, 12).Value = 2.25
$ws.Cells.Item(73, 13).Value = 3.4
$ws.Cells.Item(73, 14).Value = 3
$ws.Cells.Item(73, 15).Value = 2.6
$ws.Cells.Item(73, 16).Value = 3.3
$ws.Cells.Item(73, 17).Value = 2.625
$ws.Cells.Item(73, 19).Value = 1.85
$ws.Cells.Item(73, 20).Value = 2
$ws.Cells.Item(73, 21).Value = 2.75
$ws.Cells.Item(73, 22).Value = 1.975
$ws.Cells.Item(73, 23).Value = 1.875
$ws.Cells.Item(73, 24).Value = -1
$ws.Cells.Item(73, 26).Value = 1.625
$ws.Cells.Item(73, 27).Value = -1
$ws.Cells.Item(73, 28).Value = 1
$ws.Cells.Item(73, 29).Value = 0.4875
$ws.Cells.Item(73, 30).Value = -0.5

# Row 74
$ws.Cells.Item(74, 2).Value = 6803011
$ws.Cells.Item(74, 5).Value = 'MFK Karvina'
$ws.Cells.Item(74, 6).Value = 'Slovan Liberec'
$ws.Cells.Item(74, 7).Value = 5
$ws.Cells.Item(74, 8).Value = 2
$ws.Cells.Item(74, 9).Value = 3
$ws.Cells.Item(74, 10).Value = 1
$ws.Cells.Item(74, 11).Value = 'H'
$ws.Cells.Item(74, 12).Value = 2.7
$ws.Cells.Item(74, 13).Value = 3.25
$ws.Cells.Item(74, 14).Value = 2.5
$ws.Cells.Item(74, 15).Value = 2.75
$ws.Cells.Item(74, 16).Value = 3.25
$ws.Cells.Item(74, 17).Value = 2.55
$ws.Cells.Item(74, 19).Value = 2
$ws.Cells.Item(74, 20).Value = 1.85
$ws.Cells.Item(74, 21).Value = 2.5
$ws.Cells.Item(74, 22).Value = 1.95
$ws.Cells.Item(74, 23).Value = 1.9
$ws.Cells.Item(74, 24).Value = 1.75
$ws.Cells.Item(74, 26).Value = -1
$ws.Cells.Item(74, 27).Value = 1
$ws.Cells.Item(74, 28).Value = -1
$ws.Cells.Item(74, 29).Value = 0.95
$ws.Cells.Item(74, 30).Value = -1

# Row 133
$ws.Cells.Item(133, 2).Value = 6803071
$ws.Cells.Item(133, 5).Value = 'Slovacko'
$ws.Cells.Item(133, 6).Value = 'Bohemians 1905'
$ws.Cells.Item(133, 8).Value = 2
$ws.Cells.Item(133, 10).Value = 2
$ws.Cells.Item(133, 12).Value = 1.666
$ws.Cells.Item(133, 13).Value = 3.6
$ws.Cells.Item(133, 14).Value = 4.5
$ws.Cells.Item(133, 15).Value = 1.666
$ws.Cells.Item(133, 16).Value = 3.5
$ws.Cells.Item(133, 17).Value = 4.75
$ws.Cells.Item(133, 18).Value = -0.75
$ws.Cells.Item(133, 19).Value = 1.875
$ws.Cells.Item(133, 20).Value = 1.975
$ws.Cells.Item(133, 21).Value = 2.25
$ws.Cells.Item(133, 24).Value = 0.6659999999999999
$ws.Cells.Item(133, 27).Value = 0.875

# Row 134
$ws.Cells.Item(134, 2).Value = 6803074
$ws.Cells.Item(134, 5).Value = 'Slovan Liberec'
$ws.Cells.Item(134, 6).Value = 'FC Trinity Zlin'
$ws.Cells.Item(134, 8).Value = 3
$ws.Cells.Item(134, 10).Value = 1
$ws.Cells.Item(134, 12).Value = 1.55
$ws.Cells.Item(134, 13).Value = 4.2
$ws.Cells.Item(134, 14).Value = 5
$ws.Cells.Item(134, 15).Value = 1.5
$ws.Cells.Item(134, 16).Value = 4.2
$ws.Cells.Item(134, 17).Value = 5.75
$ws.Cells.Item(134, 18).Value = -1
$ws.Cells.Item(134, 19).Value = 1.925
$ws.Cells.Item(134, 20).Value = 1.925
$ws.Cells.Item(134, 21).Value = 2.5
$ws.Cells.Item(134, 24).Value = 0.5
$ws.Cells.Item(134, 27).Value = 0.925

# Row 143
$ws.Cells.Item(143, 2).Value = 6803081
$ws.Cells.Item(143, 5).Value = 'FC Trinity Zlin'
$ws.Cells.Item(143, 6).Value = 'Hradec Kralove'
$ws.Cells.Item(143, 7).Value = 4
$ws.Cells.Item(143, 8).Value = 0
$ws.Cells.Item(143, 9).Value = 4
$ws.Cells.Item(143, 11).Value = 'H'
$ws.Cells.Item(143, 12).Value = 2.5
$ws.Cells.Item(143, 13).Value = 3.4
$ws.Cells.Item(143, 14).Value = 2.7
$ws.Cells.Item(143, 15).Value = 3
$ws.Cells.Item(143, 16).Value = 3.4
$ws.Cells.Item(143, 17).Value = 2.3
$ws.Cells.Item(143, 18).Value = 0.25
$ws.Cells.Item(143, 19).Value = 1.825
$ws.Cells.Item(143, 20).Value = 2.025
$ws.Cells.Item(143, 22).Value = 2.05
$ws.Cells.Item(143, 23).Value = 1.8
$ws.Cells.Item(143, 24).Value = 2
$ws.Cells.Item(143, 26).Value = -1
$ws.Cells.Item(143, 27).Value = 0.825
$ws.Cells.Item(143, 28).Value = -1
$ws.Cells.Item(143, 29).Value = 1.05
$ws.Cells.Item(143, 30).Value = -1

# Row 144
$ws.Cells.Item(144, 2).Value = 6803080
$ws.Cells.Item(144, 5).Value = 'Mlada Boleslav'
$ws.Cells.Item(144, 6).Value = 'Slovacko'
$ws.Cells.Item(144, 7).Value = 0
$ws.Cells.Item(144, 8).Value = 1
$ws.Cells.Item(144, 9).Value = 0
$ws.Cells.Item(144, 11).Value = 'A'
$ws.Cells.Item(144, 12).Value = 2.3
$ws.Cells.Item(144, 13).Value = 3.3
$ws.Cells.Item(144, 14).Value = 3.1
$ws.Cells.Item(144, 15).Value = 2.5
$ws.Cells.Item(144, 16).Value = 3.3
$ws.Cells.Item(144, 17).Value = 2.8
$ws.Cells.Item(144, 18).Value = 0
$ws.Cells.Item(144, 19).Value = 1.8
$ws.Cells.Item(144, 20).Value = 2.05
$ws.Cells.Item(144, 22).Value = 1.95
$ws.Cells.Item(144, 23).Value = 1.9
$ws.Cells.Item(144, 24).Value = -1
$ws.Cells.Item(144, 26).Value = 1.8
$ws.Cells.Item(144, 27).Value = -1
$ws.Cells.Item(144, 28).Value = 1.05
$ws.Cells.Item(144, 29).Value = -1
$ws.Cells.Item(144, 30).Value = 0.8999999999999999

# Row 155
$ws.Cells.Item(155, 2).Value = 6803085
$ws.Cells.Item(155, 5).Value = 'Hradec Kralove'
$ws.Cells.Item(155, 6).Value = 'Bohemians 1905'
$ws.Cells.Item(155, 8).Value = 2
$ws.Cells.Item(155, 11).Value = 'D'
$ws.Cells.Item(155, 12).Value = 2.3
$ws.Cells.Item(155, 13).Value = 3.4
$ws.Cells.Item(155, 14).Value = 2.7
$ws.Cells.Item(155, 15).Value = 2.5
$ws.Cells.Item(155, 16).Value = 3.3
$ws.Cells.Item(155, 17).Value = 2.75
$ws.Cells.Item(155, 18).Value = 0
$ws.Cells.Item(155, 22).Value = 1.8
$ws.Cells.Item(155, 23).Value = 2.05
$ws.Cells.Item(155, 24).Value = -1
$ws.Cells.Item(155, 25).Value = 2.3
$ws.Cells.Item(155, 27).Value = 0
$ws.Cells.Item(155, 28).Value = 0
$ws.Cells.Item(155, 29).Value = 0.8

# Row 156
$ws.Cells.Item(156, 2).Value = 6803089
$ws.Cells.Item(156, 5).Value = 'FK Teplice'
$ws.Cells.Item(156, 6).Value = 'FC Trinity Zlin'
$ws.Cells.Item(156, 8).Value = 1
$ws.Cells.Item(156, 11).Value = 'H'
$ws.Cells.Item(156, 12).Value = 1.85
$ws.Cells.Item(156, 13).Value = 3.5
$ws.Cells.Item(156, 14).Value = 3.8
$ws.Cells.Item(156, 15).Value = 1.8
$ws.Cells.Item(156, 16).Value = 3.6
$ws.Cells.Item(156, 17).Value = 4.333
$ws.Cells.Item(156, 18).Value = -0.5
$ws.Cells.Item(156, 22).Value = 1.9
$ws.Cells.Item(156, 23).Value = 1.95
$ws.Cells.Item(156, 24).Value = 0.8
$ws.Cells.Item(156, 25).Value = -1
$ws.Cells.Item(156, 27).Value = 0.8
$ws.Cells.Item(156, 28).Value = -1
$ws.Cells.Item(156, 29).Value = 0.8999999999999999

# Row 162
$ws.Cells.Item(162, 2).Value = 6803100
$ws.Cells.Item(162, 5).Value = 'Pardubice'
$ws.Cells.Item(162, 6).Value = 'FK Teplice'
$ws.Cells.Item(162, 8).Value = 1
$ws.Cells.Item(162, 10).Value = 0
$ws.Cells.Item(162, 11).Value = 'D'
$ws.Cells.Item(162, 12).Value = 2.3
$ws.Cells.Item(162, 14).Value = 3
$ws.Cells.Item(162, 15).Value = 2.2
$ws.Cells.Item(162, 16).Value = 3.3
$ws.Cells.Item(162, 19).Value = 1.875
$ws.Cells.Item(162, 20).Value = 1.975
$ws.Cells.Item(162, 21).Value = 2.25
$ws.Cells.Item(162, 25).Value = 2.3
$ws.Cells.Item(162, 26).Value = -1
$ws.Cells.Item(162, 27).Value = -0.5
$ws.Cells.Item(162, 28).Value = 0.4875
$ws.Cells.Item(162, 29).Value = -0.5
$ws.Cells.Item(162, 30).Value = 0.5125

# Row 163
$ws.Cells.Item(163, 2).Value = 6803094
$ws.Cells.Item(163, 5).Value = 'Sigma Olomouc'
$ws.Cells.Item(163, 6).Value = 'Hradec Kralove'
$ws.Cells.Item(163, 7).Value = 0
$ws.Cells.Item(163, 8).Value = 0
$ws.Cells.Item(163, 9).Value = 0
$ws.Cells.Item(163, 12).Value = 1.8
$ws.Cells.Item(163, 13).Value = 3.6
$ws.Cells.Item(163, 14).Value = 4.333
$ws.Cells.Item(163, 15).Value = 1.8
$ws.Cells.Item(163, 16).Value = 3.6
$ws.Cells.Item(163, 17).Value = 4.5
$ws.Cells.Item(163, 18).Value = -0.75
$ws.Cells.Item(163, 19).Value = 2
$ws.Cells.Item(163, 20).Value = 1.85
$ws.Cells.Item(163, 21).Value = 2.5
$ws.Cells.Item(163, 22).Value = 1.975
$ws.Cells.Item(163, 23).Value = 1.875
$ws.Cells.Item(163, 25).Value = 2.6
$ws.Cells.Item(163, 27).Value = -1
$ws.Cells.Item(163, 28).Value = 0.8500000000000001
$ws.Cells.Item(163, 29).Value = -1
$ws.Cells.Item(163, 30).Value = 0.875

# Row 164
$ws.Cells.Item(164, 2).Value = 6803095
$ws.Cells.Item(164, 5).Value = 'Mlada Boleslav'
$ws.Cells.Item(164, 6).Value = 'Banik Ostrava'
$ws.Cells.Item(164, 7).Value = 1
$ws.Cells.Item(164, 8).Value = 3
$ws.Cells.Item(164, 9).Value = 1
$ws.Cells.Item(164, 10).Value = 3
$ws.Cells.Item(164, 11).Value = 'A'
$ws.Cells.Item(164, 12).Value = 2.375
$ws.Cells.Item(164, 13).Value = 3.4
$ws.Cells.Item(164, 14).Value = 2.9
$ws.Cells.Item(164, 15).Value = 2.1
$ws.Cells.Item(164, 17).Value = 3.25
$ws.Cells.Item(164, 18).Value = -0.25
$ws.Cells.Item(164, 19).Value = 1.85
$ws.Cells.Item(164, 20).Value = 2
$ws.Cells.Item(164, 21).Value = 2.75
$ws.Cells.Item(164, 22).Value = 1.825
$ws.Cells.Item(164, 23).Value = 2.025
$ws.Cells.Item(164, 25).Value = -1
$ws.Cells.Item(164, 26).Value = 2.25
$ws.Cells.Item(164, 28).Value = 1
$ws.Cells.Item(164, 29).Value = 0.825
$ws.Cells.Item(164, 30).Value = -1

# Row 166
$ws.Cells.Item(166, 2).Value = 6803097
$ws.Cells.Item(166, 5).Value = 'MFK Karvina'
$ws.Cells.Item(166, 6).Value = 'Slavia Prague'
$ws.Cells.Item(166, 7).Value = 0
$ws.Cells.Item(166, 8).Value = 3
$ws.Cells.Item(166, 10).Value = 2
$ws.Cells.Item(166, 11).Value = 'A'
$ws.Cells.Item(166, 12).Value = 8.5
$ws.Cells.Item(166, 13).Value = 6
$ws.Cells.Item(166, 14).Value = 1.285
$ws.Cells.Item(166, 15).Value = 11
$ws.Cells.Item(166, 16).Value = 6
$ws.Cells.Item(166, 17).Value = 1.25
$ws.Cells.Item(166, 18).Value = 1.75
$ws.Cells.Item(166, 19).Value = 1.875
$ws.Cells.Item(166, 20).Value = 1.975
$ws.Cells.Item(166, 21).Value = 3
$ws.Cells.Item(166, 22).Value = 1.85
$ws.Cells.Item(166, 23).Value = 2
$ws.Cells.Item(166, 24).Value = -1
$ws.Cells.Item(166, 26).Value = 0.25
$ws.Cells.Item(166, 27).Value = -1
$ws.Cells.Item(166, 28).Value = 0.9750000000000001
$ws.Cells.Item(166, 29).Value = 0
$ws.Cells.Item(166, 30).Value = 0

# Row 167
$ws.Cells.Item(167, 2).Value = 6803096
$ws.Cells.Item(167, 5).Value = 'FC Trinity Zlin'
$ws.Cells.Item(167, 6).Value = 'Slovacko'
$ws.Cells.Item(167, 7).Value = 2
$ws.Cells.Item(167, 8).Value = 1
$ws.Cells.Item(167, 10).Value = 1
$ws.Cells.Item(167, 11).Value = 'H'
$ws.Cells.Item(167, 12).Value = 3.5
$ws.Cells.Item(167, 13).Value = 3.4
$ws.Cells.Item(167, 14).Value = 2.05
$ws.Cells.Item(167, 15).Value = 3.8
$ws.Cells.Item(167, 16).Value = 3.5
$ws.Cells.Item(167, 17).Value = 1.95
$ws.Cells.Item(167, 18).Value = 0.5
$ws.Cells.Item(167, 19).Value = 1.9
$ws.Cells.Item(167, 20).Value = 1.95
$ws.Cells.Item(167, 21).Value = 2.25
$ws.Cells.Item(167, 22).Value = 1.825
$ws.Cells.Item(167, 23).Value = 2.025
$ws.Cells.Item(167, 24).Value = 2.8
$ws.Cells.Item(167, 26).Value = -1
$ws.Cells.Item(167, 27).Value = 0.8999999999999999
$ws.Cells.Item(167, 28).Value = -1
$ws.Cells.Item(167, 29).Value = 0.825
$ws.Cells.Item(167, 30).Value = -1

# Row 171
$ws.Cells.Item(171, 2).Value = 6803101
$ws.Cells.Item(171, 5).Value = 'Hradec Kralove'
$ws.Cells.Item(171, 6).Value = 'Mlada Boleslav'
$ws.Cells.Item(171, 7).Value = 0
$ws.Cells.Item(171, 8).Value = 0
$ws.Cells.Item(171, 9).Value = 0
$ws.Cells.Item(171, 12).Value = 2.7
$ws.Cells.Item(171, 13).Value = 3.2
$ws.Cells.Item(171, 14).Value = 2.55
$ws.Cells.Item(171, 15).Value = 2.9
$ws.Cells.Item(171, 16).Value = 3.2
$ws.Cells.Item(171, 17).Value = 2.4
$ws.Cells.Item(171, 18).Value = 0
$ws.Cells.Item(171, 19).Value = 2.1
$ws.Cells.Item(171, 20).Value = 1.775
$ws.Cells.Item(171, 21).Value = 2.5
$ws.Cells.Item(171, 22).Value = 1.95
$ws.Cells.Item(171, 23).Value = 1.9
$ws.Cells.Item(171, 25).Value = 2.2
$ws.Cells.Item(171, 27).Value = 0
$ws.Cells.Item(171, 28).Value = 0
$ws.Cells.Item(171, 30).Value = 0.8999999999999999

# Row 172
$ws.Cells.Item(172, 2).Value = 6803104
$ws.Cells.Item(172, 5).Value = 'Banik Ostrava'
$ws.Cells.Item(172, 6).Value = 'Bohemians 1905'
$ws.Cells.Item(172, 7).Value = 1
$ws.Cells.Item(172, 8).Value = 1
$ws.Cells.Item(172, 9).Value = 1
$ws.Cells.Item(172, 10).Value = 0
$ws.Cells.Item(172, 12).Value = 1.8
$ws.Cells.Item(172, 13).Value = 3.75
$ws.Cells.Item(172, 14).Value = 4
$ws.Cells.Item(172, 15).Value = 1.6
$ws.Cells.Item(172, 16).Value = 4
$ws.Cells.Item(172, 17).Value = 4.75
$ws.Cells.Item(172, 19).Value = 1.8
$ws.Cells.Item(172, 20).Value = 2.05
$ws.Cells.Item(172, 21).Value = 2.75
$ws.Cells.Item(172, 22).Value = 1.925
$ws.Cells.Item(172, 23).Value = 1.925
$ws.Cells.Item(172, 25).Value = 3
$ws.Cells.Item(172, 28).Value = 1.05
$ws.Cells.Item(172, 29).Value = -1
$ws.Cells.Item(172, 30).Value = 0.925

# Row 173
$ws.Cells.Item(173, 2).Value = 6803105
$ws.Cells.Item(173, 5).Value = 'FK Teplice'
$ws.Cells.Item(173, 6).Value = 'MFK Karvina'
$ws.Cells.Item(173, 7).Value = 2
$ws.Cells.Item(173, 8).Value = 2
$ws.Cells.Item(173, 10).Value = 1
$ws.Cells.Item(173, 12).Value = 1.727
$ws.Cells.Item(173, 13).Value = 3.6
$ws.Cells.Item(173, 14).Value = 4.5
$ws.Cells.Item(173, 15).Value = 1.65
$ws.Cells.Item(173, 16).Value = 3.6
$ws.Cells.Item(173, 17).Value = 5
$ws.Cells.Item(173, 18).Value = -0.75
$ws.Cells.Item(173, 19).Value = 1.875
$ws.Cells.Item(173, 20).Value = 1.975
$ws.Cells.Item(173, 22).Value = 1.875
$ws.Cells.Item(173, 23).Value = 1.975
$ws.Cells.Item(173, 25).Value = 2.6
$ws.Cells.Item(173, 27).Value = -1
$ws.Cells.Item(173, 28).Value = 0.9750000000000001
$ws.Cells.Item(173, 29).Value = 0.875
$ws.Cells.Item(173, 30).Value = -1

# Row 207
$ws.Cells.Item(207, 2).Value = 6804088
$ws.Cells.Item(207, 5).Value = 'Hradec Kralove'
$ws.Cells.Item(207, 6).Value = 'MFK Karvina'
$ws.Cells.Item(207, 7).Value = 2
$ws.Cells.Item(207, 9).Value = 1
$ws.Cells.Item(207, 10).Value = 1
$ws.Cells.Item(207, 12).Value = 1.615
$ws.Cells.Item(207, 14).Value = 5.5
$ws.Cells.Item(207, 15).Value = 1.7
$ws.Cells.Item(207, 16).Value = 3.8
$ws.Cells.Item(207, 17).Value = 5
$ws.Cells.Item(207, 18).Value = -0.75
$ws.Cells.Item(207, 19).Value = 1.9
$ws.Cells.Item(207, 20).Value = 1.95
$ws.Cells.Item(207, 21).Value = 2.5
$ws.Cells.Item(207, 24).Value = 0.7
$ws.Cells.Item(207, 27).Value = 0.45
$ws.Cells.Item(207, 28).Value = -0.5

# Row 208
$ws.Cells.Item(208, 2).Value = 6803141
$ws.Cells.Item(208, 5).Value = 'Banik Ostrava'
$ws.Cells.Item(208, 6).Value = 'FK Teplice'
$ws.Cells.Item(208, 7).Value = 4
$ws.Cells.Item(208, 9).Value = 2
$ws.Cells.Item(208, 10).Value = 0
$ws.Cells.Item(208, 12).Value = 1.6
$ws.Cells.Item(208, 14).Value = 5.25
$ws.Cells.Item(208, 15).Value = 1.533
$ws.Cells.Item(208, 16).Value = 4.2
$ws.Cells.Item(208, 17).Value = 5.75
$ws.Cells.Item(208, 18).Value = -1
$ws.Cells.Item(208, 19).Value = 1.875
$ws.Cells.Item(208, 20).Value = 1.975
$ws.Cells.Item(208, 21).Value = 2.75
$ws.Cells.Item(208, 24).Value = 0.5329999999999999
$ws.Cells.Item(208, 27).Value = 0.875
$ws.Cells.Item(208, 28).Value = -1

# Row 218
$ws.Cells.Item(218, 2).Value = 6803151
$ws.Cells.Item(218, 5).Value = 'Bohemians 1905'
$ws.Cells.Item(218, 6).Value = 'Sparta Prague'
$ws.Cells.Item(218, 8).Value = 3
$ws.Cells.Item(218, 9).Value = 1
$ws.Cells.Item(218, 10).Value = 1
$ws.Cells.Item(218, 12).Value = 5.25
$ws.Cells.Item(218, 13).Value = 4.333
$ws.Cells.Item(218, 14).Value = 1.571
$ws.Cells.Item(218, 15).Value = 5.25
$ws.Cells.Item(218, 16).Value = 4.333
$ws.Cells.Item(218, 17).Value = 1.571
$ws.Cells.Item(218, 18).Value = 1
$ws.Cells.Item(218, 19).Value = 1.8
$ws.Cells.Item(218, 20).Value = 2.05
$ws.Cells.Item(218, 22).Value = 1.95
$ws.Cells.Item(218, 23).Value = 1.9
$ws.Cells.Item(218, 26).Value = 0.571
$ws.Cells.Item(218, 28).Value = 1.05
$ws.Cells.Item(218, 29).Value = 0.95
$ws.Cells.Item(218, 30).Value = -1

# Row 219
$ws.Cells.Item(219, 2).Value = 6851033
$ws.Cells.Item(219, 5).Value = 'Mlada Boleslav'
$ws.Cells.Item(219, 6).Value = 'FK Teplice'
$ws.Cells.Item(219, 7).Value = 1
$ws.Cells.Item(219, 9).Value = 0
$ws.Cells.Item(219, 10).Value = 0
$ws.Cells.Item(219, 11).Value = 'A'
$ws.Cells.Item(219, 12).Value = 1.75
$ws.Cells.Item(219, 13).Value = 3.75
$ws.Cells.Item(219, 14).Value = 4.5
$ws.Cells.Item(219, 15).Value = 1.909
$ws.Cells.Item(219, 16).Value = 3.6
$ws.Cells.Item(219, 17).Value = 3.75
$ws.Cells.Item(219, 18).Value = -0.5
$ws.Cells.Item(219, 19).Value = 1.975
$ws.Cells.Item(219, 20).Value = 1.875
$ws.Cells.Item(219, 22).Value = 1.9
$ws.Cells.Item(219, 23).Value = 1.95
$ws.Cells.Item(219, 24).Value = -1
$ws.Cells.Item(219, 26).Value = 2.75
$ws.Cells.Item(219, 27).Value = -1
$ws.Cells.Item(219, 28).Value = 0.875
$ws.Cells.Item(219, 29).Value = 0.45
$ws.Cells.Item(219, 30).Value = -0.5

# Row 220
$ws.Cells.Item(220, 2).Value = 6803152
$ws.Cells.Item(220, 5).Value = 'Ceske Budejovice'
$ws.Cells.Item(220, 6).Value = 'Slovan Liberec'
$ws.Cells.Item(220, 7).Value = 3
$ws.Cells.Item(220, 8).Value = 2
$ws.Cells.Item(220, 11).Value = 'H'
$ws.Cells.Item(220, 12).Value = 3.1
$ws.Cells.Item(220, 13).Value = 3.5
$ws.Cells.Item(220, 14).Value = 2.2
$ws.Cells.Item(220, 15).Value = 3.4
$ws.Cells.Item(220, 16).Value = 3.5
$ws.Cells.Item(220, 17).Value = 2.05
$ws.Cells.Item(220, 18).Value = 0.25
$ws.Cells.Item(220, 19).Value = 2.05
$ws.Cells.Item(220, 20).Value = 1.8
$ws.Cells.Item(220, 24).Value = 2.4
$ws.Cells.Item(220, 26).Value = -1
$ws.Cells.Item(220, 27).Value = 1.05
$ws.Cells.Item(220, 28).Value = -1

# Row 240
$ws.Cells.Item(240, 2).Value = 6851035
$ws.Cells.Item(240, 5).Value = 'Ceske Budejovice'
$ws.Cells.Item(240, 6).Value = 'Slovacko'
$ws.Cells.Item(240, 7).Value = 2
$ws.Cells.Item(240, 8).Value = 2
$ws.Cells.Item(240, 10).Value = 1
$ws.Cells.Item(240, 12).Value = 2.875
$ws.Cells.Item(240, 13).Value = 3.1
$ws.Cells.Item(240, 14).Value = 2.25
$ws.Cells.Item(240, 15).Value = 2.7
$ws.Cells.Item(240, 16).Value = 3.1
$ws.Cells.Item(240, 17).Value = 2.55
$ws.Cells.Item(240, 18).Value = 0
$ws.Cells.Item(240, 19).Value = 1.975
$ws.Cells.Item(240, 20).Value = 1.875
$ws.Cells.Item(240, 25).Value = 2.1
$ws.Cells.Item(240, 27).Value = 0
$ws.Cells.Item(240, 28).Value = 0
$ws.Cells.Item(240, 29).Value = 1.025
$ws.Cells.Item(240, 30).Value = -1

# Row 241
$ws.Cells.Item(241, 2).Value = 6803156
$ws.Cells.Item(241, 5).Value = 'FK Jablonec'
$ws.Cells.Item(241, 6).Value = 'Mlada Boleslav'
$ws.Cells.Item(241, 7).Value = 1
$ws.Cells.Item(241, 8).Value = 1
$ws.Cells.Item(241, 10).Value = 0
$ws.Cells.Item(241, 12).Value = 2.3
$ws.Cells.Item(241, 13).Value = 3.2
$ws.Cells.Item(241, 14).Value = 2.75
$ws.Cells.Item(241, 15).Value = 2.25
$ws.Cells.Item(241, 16).Value = 3.4
$ws.Cells.Item(241, 17).Value = 3
$ws.Cells.Item(241, 18).Value = -0.25
$ws.Cells.Item(241, 19).Value = 2.025
$ws.Cells.Item(241, 20).Value = 1.825
$ws.Cells.Item(241, 25).Value = 2.4
$ws.Cells.Item(241, 27).Value = -0.5
$ws.Cells.Item(241, 28).Value = 0.4125
$ws.Cells.Item(241, 29).Value = -1
$ws.Cells.Item(241, 30).Value = 0.825

# Row 247
$ws.Cells.Item(247, 2).Value = 8157943
$ws.Cells.Item(247, 5).Value = 'Hradec Kralove'
$ws.Cells.Item(247, 6).Value = 'Sigma Olomouc'
$ws.Cells.Item(247, 7).Value = 3
$ws.Cells.Item(247, 8).Value = 1
$ws.Cells.Item(247, 9).Value = 3
$ws.Cells.Item(247, 10).Value = 1
$ws.Cells.Item(247, 12).Value = 2
$ws.Cells.Item(247, 13).Value = 3.25
$ws.Cells.Item(247, 14).Value = 3.5
$ws.Cells.Item(247, 15).Value = 1.95
$ws.Cells.Item(247, 16).Value = 3.3
$ws.Cells.Item(247, 17).Value = 3.6
$ws.Cells.Item(247, 18).Value = -0.5
$ws.Cells.Item(247, 19).Value = 1.975
$ws.Cells.Item(247, 20).Value = 1.875
$ws.Cells.Item(247, 22).Value = 1.95
$ws.Cells.Item(247, 23).Value = 1.9
$ws.Cells.Item(247, 24).Value = 0.95
$ws.Cells.Item(247, 27).Value = 0.9750000000000001
$ws.Cells.Item(247, 29).Value = 0.95
$ws.Cells.Item(247, 30).Value = -1

# Row 248
$ws.Cells.Item(248, 2).Value = 8157942
$ws.Cells.Item(248, 5).Value = 'FK Teplice'
$ws.Cells.Item(248, 6).Value = 'Slovan Liberec'
$ws.Cells.Item(248, 7).Value = 2
$ws.Cells.Item(248, 8).Value = 0
$ws.Cells.Item(248, 9).Value = 0
$ws.Cells.Item(248, 10).Value = 0
$ws.Cells.Item(248, 12).Value = 2.5
$ws.Cells.Item(248, 13).Value = 3.2
$ws.Cells.Item(248, 14).Value = 2.625
$ws.Cells.Item(248, 15).Value = 2.3
$ws.Cells.Item(248, 16).Value = 3.25
$ws.Cells.Item(248, 17).Value = 2.875
$ws.Cells.Item(248, 18).Value = -0.25
$ws.Cells.Item(248, 19).Value = 2.05
$ws.Cells.Item(248, 20).Value = 1.8
$ws.Cells.Item(248, 22).Value = 2
$ws.Cells.Item(248, 23).Value = 1.85
$ws.Cells.Item(248, 24).Value = 1.3
$ws.Cells.Item(248, 27).Value = 1.05
$ws.Cells.Item(248, 29).Value = -1
$ws.Cells.Item(248, 30).Value = 0.8500000000000001

# Row 272
$ws.Cells.Item(272, 2).Value = 8157960
$ws.Cells.Item(272, 5).Value = 'MFK Karvina'
$ws.Cells.Item(272, 6).Value = 'Ceske Budejovice'
$ws.Cells.Item(272, 7).Value = 1
$ws.Cells.Item(272, 8).Value = 0
$ws.Cells.Item(272, 9).Value = 1
$ws.Cells.Item(272, 11).Value = 'H'
$ws.Cells.Item(272, 12).Value = 2
$ws.Cells.Item(272, 13).Value = 3.4
$ws.Cells.Item(272, 14).Value = 3.3
$ws.Cells.Item(272, 15).Value = 2.4
$ws.Cells.Item(272, 16).Value = 3.1
$ws.Cells.Item(272, 17).Value = 2.8
$ws.Cells.Item(272, 18).Value = 0
$ws.Cells.Item(272, 19).Value = 1.775
$ws.Cells.Item(272, 20).Value = 2.1
$ws.Cells.Item(272, 21).Value = 2.75
$ws.Cells.Item(272, 22).Value = 2
$ws.Cells.Item(272, 23).Value = 1.85
$ws.Cells.Item(272, 24).Value = 1.4
$ws.Cells.Item(272, 26).Value = -1
$ws.Cells.Item(272, 27).Value = 0.7749999999999999
$ws.Cells.Item(272, 28).Value = -1
$ws.Cells.Item(272, 30).Value = 0.8500000000000001

# Row 273
$ws.Cells.Item(273, 2).Value = 8157958
$ws.Cells.Item(273, 5).Value = 'Bohemians 1905'
$ws.Cells.Item(273, 6).Value = 'Pardubice'
$ws.Cells.Item(273, 7).Value = 0
$ws.Cells.Item(273, 8).Value = 1
$ws.Cells.Item(273, 11).Value = 'A'
$ws.Cells.Item(273, 12).Value = 1.8
$ws.Cells.Item(273, 13).Value = 3.6
$ws.Cells.Item(273, 14).Value = 3.8
$ws.Cells.Item(273, 15).Value = 1.666
$ws.Cells.Item(273, 16).Value = 4
$ws.Cells.Item(273, 17).Value = 4.1
$ws.Cells.Item(273, 18).Value = -0.75
$ws.Cells.Item(273, 19).Value = 1.85
$ws.Cells.Item(273, 20).Value = 2
$ws.Cells.Item(273, 21).Value = 3.25
$ws.Cells.Item(273, 22).Value = 2.05
$ws.Cells.Item(273, 23).Value = 1.8
$ws.Cells.Item(273, 24).Value = -1
$ws.Cells.Item(273, 26).Value = 3.1
$ws.Cells.Item(273, 27).Value = -1
$ws.Cells.Item(273, 28).Value = 1
$ws.Cells.Item(273, 30).Value = 0.8

# Row 274
$ws.Cells.Item(274, 2).Value = 8210454
$ws.Cells.Item(274, 5).Value = 'Hradec Kralove'
$ws.Cells.Item(274, 6).Value = 'FK Teplice'
$ws.Cells.Item(274, 7).Value = 2
$ws.Cells.Item(274, 9).Value = 0
$ws.Cells.Item(274, 12).Value = 1.909
$ws.Cells.Item(274, 13).Value = 3.5
$ws.Cells.Item(274, 14).Value = 3.5
$ws.Cells.Item(274, 15).Value = 1.95
$ws.Cells.Item(274, 16).Value = 3.6
$ws.Cells.Item(274, 17).Value = 3.25
$ws.Cells.Item(274, 18).Value = -0.5
$ws.Cells.Item(274, 19).Value = 2.025
$ws.Cells.Item(274, 20).Value = 1.825
$ws.Cells.Item(274, 21).Value = 2.5
$ws.Cells.Item(274, 22).Value = 1.925
$ws.Cells.Item(274, 23).Value = 1.925
$ws.Cells.Item(274, 24).Value = 0.95
$ws.Cells.Item(274, 27).Value = 1.025
$ws.Cells.Item(274, 30).Value = 0.925

# Row 275
$ws.Cells.Item(275, 2).Value = 8157941
$ws.Cells.Item(275, 5).Value = 'Banik Ostrava'
$ws.Cells.Item(275, 6).Value = 'Slovacko'
$ws.Cells.Item(275, 7).Value = 6
$ws.Cells.Item(275, 12).Value = 1.666
$ws.Cells.Item(275, 13).Value = 4
$ws.Cells.Item(275, 14).Value = 4.2
$ws.Cells.Item(275, 15).Value = 1.7
$ws.Cells.Item(275, 16).Value = 4.1
$ws.Cells.Item(275, 17).Value = 3.9
$ws.Cells.Item(275, 18).Value = -0.75
$ws.Cells.Item(275, 21).Value = 3.25
$ws.Cells.Item(275, 22).Value = 2
$ws.Cells.Item(275, 23).Value = 1.85
$ws.Cells.Item(275, 24).Value = 0.7
$ws.Cells.Item(275, 29).Value = 1
$ws.Cells.Item(275, 30).Value = -1

# Row 276
$ws.Cells.Item(276, 2).Value = 8157175
$ws.Cells.Item(276, 5).Value = 'Slavia Prague'
$ws.Cells.Item(276, 6).Value = 'Mlada Boleslav'
$ws.Cells.Item(276, 7).Value = 4
$ws.Cells.Item(276, 12).Value = 1.166
$ws.Cells.Item(276, 13).Value = 6.5
$ws.Cells.Item(276, 14).Value = 13
$ws.Cells.Item(276, 15).Value = 1.222
$ws.Cells.Item(276, 16).Value = 6
$ws.Cells.Item(276, 17).Value = 11
$ws.Cells.Item(276, 18).Value = -2
$ws.Cells.Item(276, 21).Value = 3.75
$ws.Cells.Item(276, 22).Value = 1.875
$ws.Cells.Item(276, 23).Value = 1.975
$ws.Cells.Item(276, 24).Value = 0.222
$ws.Cells.Item(276, 29).Value = 0.4375
$ws.Cells.Item(276, 30).Value = -0.5
